$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (changed) date column C was bumped from 45186 to 45188
# for every existing data row (2 .. 425).
$ws.Range("C2:C425").Value = 45188

# Row 425 picked up an explicit row-height attribute in the diff.
$ws.Rows.Item(425).RowHeight = 15

# New data row 426 appended at the bottom of the table.
$ws.Range("A426").Value = "A 43689-2023"

$ws.Range("B426").Value = 45187
$ws.Range("B426").NumberFormat = $ws.Range("B425").NumberFormat

$ws.Range("C426").Value = 45188
$ws.Range("C426").NumberFormat = $ws.Range("C425").NumberFormat

$ws.Range("D426").Value = "SKÅNE LÄN"
$ws.Range("E426").Value = "OSBY"

$ws.Range("G426").Value = 2.5
$ws.Range("H426:Q426").Value = 0

$ws.Range("R426").WrapText = $true
